# Generate Report for Handoff
# Re-stamps the "Latest Handoff Date"/"Latest Handoff Datetime" columns for every
# file that is still pending (i.e. not already "Handed back: in sync with en-US"
# or "In Translation") with the new report-generation timestamp.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D is "Latest Handoff Date" ---
$overview = $wb.Worksheets.Item("Overview")
$overviewRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $overviewRows) {
    $overview.Cells.Item($r, 4).Value = "2016-22-12 16:22:25"
}

# --- zh-cn sheet: column E is "Latest Handoff Datetime" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$langRows = @(7, 10, 11, 12, 13, 14, 15, 16)
foreach ($r in $langRows) {
    $zhcn.Cells.Item($r, 5).Value = "2016-03-12 16:22:22"
}

# --- de-de sheet: column E is "Latest Handoff Datetime" ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $langRows) {
    $dede.Cells.Item($r, 5).Value = "2016-03-12 16:22:25"
}
